$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 1 to push everything down, and make room for column B.
$ws.Rows.Item(1).Insert()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 54.140625
$ws.Columns.Item(2).ColumnWidth = 51.7109375

# Fill content in authoring order (matches shared-string insertion order)
$ws.Cells.Item(1,2).Value = "ghi chú"
$ws.Cells.Item(2,2).Value = "post hiển thị bắt buộc phải có technology không sẽ bị lỗi"
$ws.Cells.Item(1,1).Value = "Chức năng"
$ws.Cells.Item(4,2).Value = "điểm danh tiếng của bài viết tính theo lượt vow của bài post + trả lời + đề xuất+ lượt xem"
$ws.Cells.Item(3,2).Value = "huy hiệu tính theo bài viết + câu trả lời : bài viết được 1 vow trở lên 1 đồng, 4 vow 1 bạc, 10 vow vàng, 20 vow bạc kim, 30 vow ruby"
$ws.Cells.Item(5,2).Value = "chức năng các nhân vật đánh giá nhau"

# Formatting for header row
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Size = 20
$headerRange.Font.Color = 255
$headerRange.Font.Name = "Arial"
$headerRange.Interior.ThemeColor = 2
$headerRange.HorizontalAlignment = -4108
$ws.Cells.Item(1,1).VerticalAlignment = -4108
$headerRange.RowHeight = 25.5

# Wrap text + alignment for B3:B4
$ws.Range("B3:B4").WrapText = $true

# Page setup
$ws.PageSetup.Orientation = 1

$wb.Save()
